$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The sheet had two columns both labelled "cdr" by mistake (one for the
# roll-due-to-rudder derivative, one for the yaw-due-to-rudder derivative).
# Rename them to their correct, distinct coefficient names.
$ws.Range("M1").Value = "cldr"
$ws.Range("Q1").Value = "cndr"

# Updated simulation results - more accurate now, but still off.
$ws.Range("L2").Value = 0.03
$ws.Range("N2").Value = 0.1

# Restore the on-screen selection/scroll state that was active when saved.
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N3").Select()
